$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update activity descriptions in column C
$ws.Range("C7").Value = "Se creo el nuevo repositorio"
$ws.Range("C8").Value = "Se creo el proyecto y la aplicación"
$ws.Range("C9").Value = "se creo la carpeta static y agrega las imágenes, css, js"
$ws.Range("C10").Value = "Se comenzo a adaptar el proyecto a Django "
$ws.Range("C11").Value = "Se comenzó a programar las nuevas paginas de carrito y comprar"
$ws.Range("C15").Value = "Se programo la opción de Créate del CRUD en el formulario para agregar un producto"
$ws.Range("C16").Value = "Se programo la opcion de Update del CRUD en el formulario para modificar un producto"
$ws.Range("C17").Value = "Se arreglaron algunos detalles del proyecto"
$ws.Range("C18").Value = "Se realizo el ultimo commit"

# Widen column C (target stored width 66.6719; engine rounds ColumnWidth to the
# nearest 1/6 of a character, so 65.75 is the closest achievable input)
$ws.Columns.Item(3).ColumnWidth = 65.75
